$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare cell formatting (styles / row heights) for the new rows first ---
# Row 13 keeps its existing style (s=6/7), just needs a value added in column A.
$ws.Rows.Item(13).RowHeight = 43.2

# Row 14: new data row -> copy cell styles from row 12 (s=4/4/5/5/5 pattern).
$ws.Range("A12:E12").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)
$ws.Rows.Item(14).RowHeight = 43.2

# Row 15: new data row (no filename cell) -> copy cell styles from row 10 (s=4/5 pattern, no A cell).
$ws.Range("B10:E10").Copy()
$ws.Range("B15:E15").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Now fill in the values, in the same order the original author typed them ---
$ws.Range("A13").Value = 'SCRIPT/T01P01A/um2507.ssb'
$ws.Range("C14").Value = ' You did it! That\''s wonderful!'
$ws.Range("C15").Value = ' Thank you so much!'
$ws.Range("A14").Value = 'SCRIPT/T01P02A/us0104.ssb'
$ws.Range("D14").Value = ' У вас всё получилось!\nКак здорово!'
$ws.Range("D15").Value = ' Спасибо вам огромное!'
$ws.Range("E14").Value = ' Ô âàò âòæ ðïìôœéìïòû!\nËàë èäïñïâï!'
$ws.Range("E15").Value = ' Òðàòéáï âàí ïãñïíîïå!'

$ws.Range("B14").Value = 60
$ws.Range("B15").Value = 63

$ws.Range("E15").Select()
